# Fix inconsistencies in the "Inverse Property" column text (remove
# diacritics so the strings match the de-accented style used elsewhere
# in the sheet), then leave the active selection on the last data cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: CasovaVyuzitelnost -> jePouzivanVRocnimObdobi -> jeObdobimPouzitiPro
$ws.Range("D6").Value = "jeObdobimPouzitiPro"

# Row 9: MistoVyuziti -> jeMistneVyuzitelnyV -> jeMistemPouzitiPro
$ws.Range("D9").Value = "jeMistemPouzitiPro"

# Move/update the selection to D14 (last cell of the table)
$ws.Range("D14").Select() | Out-Null
